$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "[49.83440990109228, 50.07240287229322]"
$ws.Range("U2").Value = "[49.8510931545863, 50.025301190358896]"

$ws.Range("M3").Value = "[49.83118640922864, 50.13303153339768]"
$ws.Range("U3").Value = "[49.898704568608245, 50.06537575287347]"

$ws.Range("M4").Value = "[49.947491046461614, 50.26520811537957]"
$ws.Range("U4").Value = "[49.90863780043246, 50.0720571221768]"

$ws.Range("M5").Value = "[49.83390281958106, 50.132281062687404]"
$ws.Range("U5").Value = "[49.906021828148376, 50.07820704958349]"

$ws.Range("M6").Value = "[49.94069116445341, 50.22654684561526]"
$ws.Range("U6").Value = "[49.96379313189229, 50.13226112767495]"

$ws.Range("M7").Value = "[49.99560003809346, 50.25459162277086]"
$ws.Range("U7").Value = "[49.985582925195395, 50.15209189452805]"

$ws.Range("M8").Value = "[49.906662952694774, 50.16393849006451]"
$ws.Range("U8").Value = "[49.86777179917612, 50.0500161667644]"
